$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update rnaSampleNumber column (C) for rows 2-27: add 26 to each existing value
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $cell.Value2 + 26
}

# Update the selection on the active sheet to C2:C27 with active cell C2
$ws.Range("C2:C27").Select()
